# Updated symbol list on Thu Dec 15 22:45:19 UTC 2022 with GitHub Actions
# Refresh the "Price" column (D) values for the coins whose quotes moved.
# The Price column stores text (numbers-as-text), so each value is entered
# with a leading apostrophe to force text entry, matching the original
# inlineStr/text storage rather than being auto-converted to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'257.41"
$ws.Range("D3").Value  = "'22.81"
$ws.Range("D4").Value  = "'6.162"
$ws.Range("D5").Value  = "'0.06063"
$ws.Range("D9").Value  = "'0.7970"
$ws.Range("D11").Value = "'0.08062"
$ws.Range("D12").Value = "'0.03356"
$ws.Range("D13").Value = "'0.03082"
$ws.Range("D14").Value = "'0.09300"
$ws.Range("D15").Value = "'3.915"
$ws.Range("D16").Value = "'0.001696"
$ws.Range("D17").Value = "'0.04830"
$ws.Range("D18").Value = "'0.0006158"
$ws.Range("D19").Value = "'0.006214"
$ws.Range("D20").Value = "'0.001101"
$ws.Range("D21").Value = "'0.003379"
$ws.Range("D22").Value = "'0.0001504"
$ws.Range("D23").Value = "'3.685"
$ws.Range("D24").Value = "'2.264"
$ws.Range("D26").Value = "'0.1227"
$ws.Range("D27").Value = "'0.0003023"
$ws.Range("D40").Value = "'0.04578"
$ws.Range("D41").Value = "'0.007167"
$ws.Range("D42").Value = "'0.003909"
$ws.Range("D44").Value = "'0.009939"
$ws.Range("D46").Value = "'0.00005933"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D48").Value = "'0.7518"
$ws.Range("D49").Value = "'0.06661"
$ws.Range("D50").Value = "'0.00001504"
$ws.Range("D51").Value = "'0.01012"
